# Update the dSF column (F) values to reflect repulled data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -3
$ws.Range("F4").Value = 15
$ws.Range("F6").Value = -9
$ws.Range("F7").Value = 6
$ws.Range("F8").Value = 0
$ws.Range("F11").Value = -8
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = -6
$ws.Range("F15").Value = 6
$ws.Range("F17").Value = -4
